# Apply the cryptos-list refresh described in the commit:
# "Updated cryptos list on Sun Jun 25 03:12:19 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '30.731.09'
$ws.Range("E2").Value = '  -0.05%  '

# Row 3
$ws.Range("D3").Value = '1.889.97'
$ws.Range("E3").Value = '  -0.26%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.11%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.28'
$ws.Range("E5").Value = '  -2.67%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.26%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4886'
$ws.Range("E7").Value = '  -0.98%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2975'
$ws.Range("E8").Value = '  +0.62%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06741'
$ws.Range("E9").Value = '  -1.13%  '

# Row 10
$ws.Range("D10").Value = '1.883.71'
$ws.Range("E10").Value = '  -0.54%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '16.93'
$ws.Range("E11").Value = '  -3.28%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07323'
$ws.Range("E12").Value = '  +0.79%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '89.79'
$ws.Range("E13").Value = '  -1.66%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.106'
$ws.Range("E14").Value = '  +0.21%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6670'
$ws.Range("E15").Value = '  -2.49%  '

# Row 16
$ws.Range("D16").Value = '30.672.85'
$ws.Range("E16").Value = '  -0.12%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000007949'
$ws.Range("E17").Value = '  -0.98%  '

# Row 18
$ws.Range("B18").Value = 'Avalanche'
$ws.Range("C18").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.43'
$ws.Range("E18").Value = '  +0.68%  '

# Row 19
$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.002'
$ws.Range("E19").Value = '  +0.22%  '

# Row 20
$ws.Range("D20").Value = '2.132.50'
$ws.Range("E20").Value = '  -0.04%  '

# Row 21
$ws.Range("E21").Value = '  +0.19%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '208.82'
$ws.Range("E22").Value = '  +7.23%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.952'
$ws.Range("E23").Value = '  +1.77%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.220'
$ws.Range("E24").Value = '  +1.92%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.642'
$ws.Range("E25").Value = '  +3.09%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.07'
$ws.Range("E26").Value = '  +2.62%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.97'
$ws.Range("E27").Value = '  -2.36%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.873'
$ws.Range("E28").Value = '  -3.14%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.419'
$ws.Range("E29").Value = '  +1.75%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.336'
$ws.Range("E30").Value = '  -0.27%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09140'
$ws.Range("E31").Value = '  +1.17%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.027'
$ws.Range("E32").Value = '  -0.43%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05177'
$ws.Range("E33").Value = '  -0.27%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7560'
$ws.Range("E34").Value = '  +1.02%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.107'
$ws.Range("E35").Value = '  -2.33%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.696'
$ws.Range("E36").Value = '  -0.12%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01835'
$ws.Range("E37").Value = '  -2.76%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.699'
$ws.Range("E38").Value = '  +0.67%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9259'
$ws.Range("E39").Value = '  -1.04%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.097'
$ws.Range("E40").Value = '  -3.40%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4517'
$ws.Range("E41").Value = '  +1.12%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '106.60'
$ws.Range("E42").Value = '  +0.69%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.817'
$ws.Range("E43").Value = '  -0.32%  '

# Row 44
$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.883'
$ws.Range("E44").Value = '  +2.06%  '

# Row 45
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.000'
$ws.Range("E45").Value = '  +0.05%  '

# Row 46
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1364'
$ws.Range("E46").Value = '  +1.53%  '

# Row 47
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '65.18'
$ws.Range("E47").Value = '  +12.39%  '

# Row 48
$ws.Range("B48").Value = 'Elrond'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '35.15'
$ws.Range("E48").Value = '  +4.53%  '

# Row 49
$ws.Range("B49").Value = 'Decentraland'
$ws.Range("C49").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4098'
$ws.Range("E49").Value = '  +2.97%  '

# Row 50
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05892'
$ws.Range("E50").Value = '  +0.61%  '

# Row 51
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.847'
$ws.Range("E51").Value = '  +1.19%  '
